# Updates cryptos list figures (price + 1h volume change) to match the
# latest scrape, per commit "Updated cryptos list ... with GitHub Actions".
#
# Some Price (column D) values look like plain decimals (e.g. "208.55",
# "1.00", "0.0972"); if assigned directly, Excel/COM would parse them as
# numbers and lose the original text formatting (and values like "1.00"
# would collapse to "1"). To keep them as literal text we briefly force
# the cell to Text format before assigning, then restore the (unstyled)
# "Normal" style so no stray formatting is left behind. Values that are
# not valid numbers as text (e.g. "27.582.69", with two dots) do not need
# this and are assigned directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.582.69'
$ws.Range("E2").Value = '  -1.03%  '

$ws.Range("D3").Value = '1.597.18'
$ws.Range("E3").Value = '  -1.92%  '

$ws.Range("E4").Value = '  +0.51%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '208.55'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.26%  '

$ws.Range("E6").Value = '  -3.39%  '

$ws.Range("E7").Value = '  +0.58%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.38'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.07%  '

$ws.Range("E9").Value = '  -1.73%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0592'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.14%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0865'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.71%  '

$ws.Range("D12").Value = '1.825.14'
$ws.Range("E12").Value = '  -1.79%  '

$ws.Range("D13").Value = '1.610.59'
$ws.Range("E13").Value = '  -1.04%  '

$ws.Range("E14").Value = '  -3.80%  '

$ws.Range("E15").Value = '  -4.22%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.53'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.61%  '

$ws.Range("D17").Value = '27.607.41'
$ws.Range("E17").Value = '  -0.90%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '217.63'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -5.05%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.40'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.46%  '

$ws.Range("E20").Value = '  -3.43%  '

$ws.Range("E21").Value = '  +0.53%  '

$ws.Range("E22").Value = '  -3.10%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.70'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.78%  '

$ws.Range("E24").Value = '  -1.94%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.60'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.31%  '

# Row 26/27 swapped position (BinanceUSD now ranks above Cosmos) plus
# updated price/volume figures.
$ws.Range("B26").Value = 'BinanceUSD'
$ws.Range("C26").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.53%  '

$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.74'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.11%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.09'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.78%  '

$ws.Range("E29").Value = '  -3.77%  '

$ws.Range("E30").Value = '  -1.17%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0469'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.43%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.26'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.14%  '

$ws.Range("D33").Value = '1.370.73'
$ws.Range("E33").Value = '  -1.58%  '

$ws.Range("E34").Value = '  -4.46%  '

$ws.Range("E35").Value = '  -3.10%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.975'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.17%  '

$ws.Range("E37").Value = '  -0.76%  '

$ws.Range("E38").Value = '  -2.34%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.540'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.06%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.815'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.25%  '

$ws.Range("E41").Value = '  +0.50%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.976'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.67%  '

$ws.Range("E43").Value = '  -0.70%  '

$ws.Range("E44").Value = '  -3.62%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '64.15'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.43%  '

$ws.Range("D46").Value = '1.735.38'
$ws.Range("E46").Value = '  -1.79%  '

$ws.Range("E47").Value = '  -1.83%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.96'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.10%  '

$ws.Range("E49").Value = '  -3.26%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0972'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.40%  '

$ws.Range("E51").Value = '  -0.83%  '
